$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks first - they'll be rebuilt (with shifted
# targets) after the NOC rows are removed and the data shifts up.
$ws.Hyperlinks.Delete()

# Remove the "NOC 1" / "NOC 2" / "NOC 3" rows (rows 108-110). This shifts
# FrontLiner / ICT Supervisor / BTMR / Contractor rows up by three, leaving
# the sheet with rows 1-119 instead of 1-122.
$ws.Rows("108:110").Delete()

# Replace the numeric "group id" in column G with a descriptive role code
# for every staff row (102-119).
$ws.Range("G102").Value = "BTMR"
$ws.Range("G103").Value = "BTMR"
$ws.Range("G104").Value = "BTMR"
$ws.Range("G105").Value = "JIM"
$ws.Range("G106").Value = "JIM"
$ws.Range("G107").Value = "JIM"
$ws.Range("G108").Value = "FRONTLINER"
$ws.Range("G109").Value = "FRONTLINER"
$ws.Range("G110").Value = "FRONTLINER"
$ws.Range("G111").Value = "ICT_SV"
$ws.Range("G112").Value = "ICT_SV"
$ws.Range("G113").Value = "ICT_SV"
$ws.Range("G114").Value = "BTMR"
$ws.Range("G115").Value = "BTMR"
$ws.Range("G116").Value = "BTMR"
$ws.Range("G117").Value = "CONTRACTOR"
$ws.Range("G118").Value = "CONTRACTOR"
$ws.Range("G119").Value = "CONTRACTOR"

# Re-create the mailto hyperlinks for column D, rows 105-119 (the emails
# already sit in those cells - Hyperlinks.Add just wires up the link).
$ws.Hyperlinks.Add($ws.Range("D105"), "mailto:jim1@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D106"), "mailto:jim2@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D107"), "mailto:jim3@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D108"), "mailto:frontliner1@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D109"), "mailto:frontliner2@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D110"), "mailto:frontliner3@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D111"), "mailto:ict1@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D112"), "mailto:ict2@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D113"), "mailto:ict3@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D114"), "mailto:btmr1@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D115"), "mailto:btmr2@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D116"), "mailto:btmr3@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D117"), "mailto:contractor1@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D118"), "mailto:contractor2@heitech.com.my")
$ws.Hyperlinks.Add($ws.Range("D119"), "mailto:contractor3@heitech.com.my")

# Match the final cursor/scroll position left behind in the saved file.
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F125").Select()
